# Apply "add sample data raw in bulk upload account details" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text update (L1 -> "Bookkeeping Monthly/Clean up catch up (Specify Period)") ---
$ws.Range("L1").Value = "Bookkeeping Monthly/Clean up catch up (Specify Period)"

# --- Column K width change (target stored width 30.53; engine quantizes in
# ~1/6 character steps, so 29.6 is the closest input that rounds to 30.5,
# the nearest achievable stored width) ---
$ws.Columns.Item(11).ColumnWidth = 29.6

# --- New sample data row (row 2) ---
$ws.Range("A2").Value = "AAA"
$ws.Range("B2").Value = "XYZ Company"
$ws.Range("C2").Value = "Residential Construction"
$ws.Range("D2").Value = "Real Estate"
$ws.Range("E2").Value = "Scorp"
$ws.Range("F2").Value = "QBO"
$ws.Range("G2").Value = "Dropbox"
$ws.Range("H2").Value = "Yes or No"
$ws.Range("I2").Value = "Cash or Accruals"
$ws.Range("J2").Value = 20

$ws.Range("K2").Value = "1. Categorize (Put checks in Misc)" + [char]10 + "2. Reconcile Bank Account (1 bank)"
$ws.Range("K2").WrapText = $true

$ws.Range("L2").Value = "catch up for 2021"

$ws.Range("M2").Value = 45451
$ws.Range("M2").NumberFormat = "dd/mm/yy"

$ws.Range("N2").Value = "Bank account transactions 170" + [char]10 + "Journal entries 2" + [char]10 + "Bank deposits 12" + [char]10 + "Expenses 158" + [char]9
$ws.Range("N2").WrapText = $true

$ws.Range("O2").Value = "XYZ"

# --- Row height for the new data row ---
$ws.Rows.Item(2).RowHeight = 55.2

# --- Selection / view state ---
$ws.Range("N9").Select() | Out-Null
